$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update room header names to include room numbers ---
$ws.Range("A2").Value = "RH1 Kitchen (Room1)"
$ws.Range("A10").Value = "RH1 Bath (Room2)"
$ws.Range("A15").Value = "RH1 Master Bath (Room3)"
$ws.Range("A20").Value = "RH2 Kitchen (Room4)"
$ws.Range("A30").Value = "RH2 Bath (Room5)"
$ws.Range("A35").Value = "RH2 Master Bath  (Room6)"
$ws.Range("A39").Value = "RH3 Kitchen (Room7)"
$ws.Range("A49").Value = "RH3 Bath (Room8)"
$ws.Range("A54").Value = "RH3 Master Bath (Room9)"
$ws.Range("A58").Value = "RH4 Kitchen (Room10)"
$ws.Range("A68").Value = "RH4 Bath (Room11)"
$ws.Range("A73").Value = "RH4 Master Bath (Room12)"
$ws.Range("A77").Value = "RH5 Kitchen (Room13)"
$ws.Range("A87").Value = "RH5 Bath (Room14)"
$ws.Range("A92").Value = "RH5 Master Bath (Room15)"
$ws.Range("A96").Value = "RH6 Kitchen (Room16)"
$ws.Range("A106").Value = "RH6 Bath (Room17)"
$ws.Range("A111").Value = "RH6 Master Bath (Room18)"
$ws.Range("A115").Value = "RH7 Kitchen (Room19)"
$ws.Range("A125").Value = "RH7 Bath (Room20)"
$ws.Range("A130").Value = "RH7 Master Bath (Room21)"

# --- Font size change (11 -> 10) applies to all cells using the shared "sz 11, no other attrs"
# font (fontId 4): both the A/C body-row cells and the B (CabNo) body-row cells.
# Indent change (2 -> 1) applies to the A/C body-row cells and the "Product Name"/"Notes"
# header-row cells (A/C), but NOT the centered CabNo / centered header cells.

# A/C body rows: font size 10 AND indent level 1
$bodyAC = @(
    "A4", "C4", "A5", "C5", "A6", "C6", "A7", "C7", "A8", "C8",
    "A12", "C12", "A13", "C13", "A17", "C17", "A18", "C18", "A22", "C22",
    "A23", "C23", "A24", "C24", "A25", "C25", "A26", "C26", "A27", "C27",
    "A28", "C28", "A32", "C32", "A33", "C33", "A37", "C37", "A41", "C41",
    "A42", "C42", "A43", "C43", "A44", "C44", "A45", "C45", "A46", "C46",
    "A47", "C47", "A51", "C51", "A52", "C52", "A56", "C56", "A60", "C60",
    "A61", "C61", "A62", "C62", "A63", "C63", "A64", "C64", "A65", "C65",
    "A66", "C66", "A70", "C70", "A71", "C71", "A75", "C75", "A79", "C79",
    "A80", "C80", "A81", "C81", "A82", "C82", "A83", "C83", "A84", "C84",
    "A85", "C85", "A89", "C89", "A90", "C90", "A94", "C94", "A98", "C98",
    "A99", "C99", "A100", "C100", "A101", "C101", "A102", "C102", "A103", "C103",
    "A104", "C104", "A108", "C108", "A109", "C109", "A113", "C113", "A117", "C117",
    "A118", "C118", "A119", "C119", "A120", "C120", "A121", "C121", "A122", "C122",
    "A123", "C123", "A127", "C127", "A128", "C128", "A132", "C132"
)
foreach ($c in $bodyAC) {
    $rng = $ws.Range($c)
    $rng.Font.Size = 10
    $rng.IndentLevel = 1
}

# B body rows (CabNo values): font size 10 only (indent/alignment unchanged)
$bodyB = @(
    "B4", "B5", "B6", "B7", "B8", "B12", "B13", "B17", "B18", "B22",
    "B23", "B24", "B25", "B26", "B27", "B28", "B32", "B33", "B37", "B41",
    "B42", "B43", "B44", "B45", "B46", "B47", "B51", "B52", "B56", "B60",
    "B61", "B62", "B63", "B64", "B65", "B66", "B70", "B71", "B75", "B79",
    "B80", "B81", "B82", "B83", "B84", "B85", "B89", "B90", "B94", "B98",
    "B99", "B100", "B101", "B102", "B103", "B104", "B108", "B109", "B113", "B117",
    "B118", "B119", "B120", "B121", "B122", "B123", "B127", "B128", "B132"
)
foreach ($c in $bodyB) {
    $ws.Range($c).Font.Size = 10
}

# "Product Name"/"Notes" header rows (A/C): indent level 1 only (font size unchanged)
$headerAC = @(
    "A3", "C3", "A11", "C11", "A16", "C16", "A21", "C21", "A31", "C31",
    "A36", "C36", "A40", "C40", "A50", "C50", "A55", "C55", "A59", "C59",
    "A69", "C69", "A74", "C74", "A78", "C78", "A88", "C88", "A93", "C93",
    "A97", "C97", "A107", "C107", "A112", "C112", "A116", "C116", "A126", "C126",
    "A131", "C131"
)
foreach ($c in $headerAC) {
    $ws.Range($c).IndentLevel = 1
}
